# Backup QR Scanner data - 14/08/2025, 8:32:38 AM
#
# Renames the "Scanner" worksheet to "Anatomy" and removes the last
# (4th) logged scan row, which shrinks the used range from A1:F4 to
# A1:F3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet: "Scanner" -> "Anatomy"
$ws.Name = "Anatomy"

# Remove the 4th row of scan data (student 333333, logged 08:31:41),
# shifting the remaining rows up and shrinking the sheet's used range.
$ws.Rows(4).Delete()

Write-Output "Renamed sheet to 'Anatomy' and removed row 4."
